$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so numeric-looking strings
# (e.g. "535.48") are not auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.013.37"
$ws.Range("D3").Value = "2.505.94"
$ws.Range("D5").Value = "535.48"
$ws.Range("D6").Value = "136.97"
$ws.Range("D8").Value = "0.561"
$ws.Range("D9").Value = "2.535.45"
$ws.Range("D13").Value = "0.348"
$ws.Range("D14").Value = "2.954.29"
$ws.Range("D15").Value = "23.17"
$ws.Range("D16").Value = "58.959.98"
$ws.Range("D18").Value = "2.521.63"
$ws.Range("D19").Value = "11.06"
$ws.Range("D20").Value = "4.26"
$ws.Range("D21").Value = "326.06"
$ws.Range("D23").Value = "5.87"
$ws.Range("D24").Value = "63.41"
$ws.Range("D29").Value = "6.82"
$ws.Range("D30").Value = "0.0₃0778"
$ws.Range("D31").Value = "1.77"
$ws.Range("D32").Value = "166.62"
$ws.Range("D33").Value = "1.14"
$ws.Range("D34").Value = "0.997"
$ws.Range("D35").Value = "1.40"
$ws.Range("D36").Value = "18.48"
$ws.Range("D37").Value = "4.12"
$ws.Range("D39").Value = "36.80"
$ws.Range("D40").Value = "0.826"
$ws.Range("D42").Value = "5.24"
$ws.Range("D43").Value = "278.50"
$ws.Range("D45").Value = "0.605"
$ws.Range("D46").Value = "10.85"
$ws.Range("D47").Value = "125.59"
$ws.Range("D48").Value = "0.0927"
$ws.Range("D49").Value = "0.0512"
$ws.Range("D51").Value = "17.57"

$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E6").Value = "  -2.03%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("E21").Value = "  +1.92%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("E37").Value = "  -2.82%  "
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E40").Value = "  +2.82%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("E50").Value = "  -0.15%  "

# Restore the original (default) cell formatting now that the text values are set.
$ws.Range("D2:D51").ClearFormats()
